$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 90
$ws.Range("E5").Value = 123
$ws.Range("E7").Value = 30
$ws.Range("E10").Value = 439
$ws.Range("F10").Value = 216
$ws.Range("H10").Value = 216
$ws.Range("E11").Value = 300
$ws.Range("F11").Value = 165
$ws.Range("H11").Value = 165
$ws.Range("E12").Value = 435
$ws.Range("E14").Value = 110
$ws.Range("E16").Value = 181
$ws.Range("F16").Value = 91
$ws.Range("H16").Value = 91
$ws.Range("E17").Value = 86
$ws.Range("E21").Value = 130
$ws.Range("E22").Value = 151
$ws.Range("E23").Value = 181
$ws.Range("E24").Value = 185
$ws.Range("F25").Value = 104
$ws.Range("H25").Value = 104
$ws.Range("E26").Value = 134
$ws.Range("F26").Value = 78
$ws.Range("H26").Value = 78
$ws.Range("E28").Value = 180
$ws.Range("E29").Value = 154
$ws.Range("F29").Value = 81
$ws.Range("H29").Value = 81
$ws.Range("E30").Value = 189
$ws.Range("E32").Value = 168
$ws.Range("E33").Value = 259
$ws.Range("F33").Value = 129
$ws.Range("H33").Value = 129
$ws.Range("E34").Value = 196
$ws.Range("F34").Value = 117
$ws.Range("H34").Value = 117
$ws.Range("E35").Value = 131
$ws.Range("E36").Value = 62
$ws.Range("E37").Value = 141
$ws.Range("F37").Value = 67
$ws.Range("H37").Value = 67
$ws.Range("E40").Value = 236
$ws.Range("E41").Value = 356
$ws.Range("F41").Value = 163
$ws.Range("H41").Value = 163
$ws.Range("E42").Value = 326
$ws.Range("F42").Value = 170
$ws.Range("H42").Value = 170
$ws.Range("E43").Value = 107
$ws.Range("F43").Value = 58
$ws.Range("H43").Value = 58
$ws.Range("E45").Value = 126
$ws.Range("F45").Value = 62
$ws.Range("H45").Value = 62
$ws.Range("E46").Value = 280
$ws.Range("F46").Value = 150
$ws.Range("H46").Value = 150
$ws.Range("E47").Value = 394
$ws.Range("E49").Value = 263
$ws.Range("E50").Value = 228
$ws.Range("F50").Value = 98
$ws.Range("H50").Value = 98
$ws.Range("E52").Value = 24

Write-Output "Done updating Inscritos/Pagos/Homologadas cells"
